# Applies the "Tolto parametro a / beta dal main per testaggio" re-run of the
# swap_intra scheduling simulation: rows 5-7 now reflect a different ordering
# of jobs on machine "BIMEC 5" (job 251752 scheduled first, then 251218, then
# 251895), with updated setup/processing timestamps and lateness (column R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (was commessa 251218, now commessa 251752) ---
$ws.Cells.Item(5, 1).Value  = 251752                                    # A5 commessa
$ws.Cells.Item(5, 3).Value  = 15                                        # C5 minuti setup
$ws.Cells.Item(5, 4).Value  = 0                                         # D5 minuti processamento
$ws.Cells.Item(5, 6).Value  = "2025-06-04 13:07:16"                     # F5 fine setup
$ws.Cells.Item(5, 7).Value  = "2025-06-04 13:07:16"                     # G5 inizio lavorazione
$ws.Cells.Item(5, 8).Value  = "2025-06-04 13:07:16"                     # H5 fine lavorazione
$ws.Cells.Item(5, 9).Value  = 0                                         # I5 mt da tagliare
$ws.Cells.Item(5, 11).Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R6 ;R9"  # K5 macchine compatibili
$ws.Cells.Item(5, 12).Value = 3                                         # L5 numero coltelli
$ws.Cells.Item(5, 14).Value = 39846                                     # N5 veicolo
$ws.Cells.Item(5, 16).Value = 39846                                     # P5
$ws.Cells.Item(5, 17).Value = "2025-05-20 00:00:00"                     # Q5
$ws.Cells.Item(5, 18).Value = -0.5467233959259259                       # R5

# --- Row 6 (was commessa 251895, now commessa 251218) ---
$ws.Cells.Item(6, 1).Value  = 251218                                    # A6 commessa
$ws.Cells.Item(6, 3).Value  = 21                                        # C6 minuti setup
$ws.Cells.Item(6, 4).Value  = 96.90140845070422                         # D6 minuti processamento
$ws.Cells.Item(6, 5).Value  = "2025-06-04 13:07:16"                     # E6 inizio setup
$ws.Cells.Item(6, 6).Value  = "2025-06-04 13:28:16"                     # F6 fine setup
$ws.Cells.Item(6, 7).Value  = "2025-06-04 13:28:16"                     # G6 inizio lavorazione
$ws.Cells.Item(6, 8).Value  = "2025-06-05 07:05:10"                     # H6 fine lavorazione
$ws.Cells.Item(6, 9).Value  = 6880                                      # I6 mt da tagliare
$ws.Cells.Item(6, 11).Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9"      # K6 macchine compatibili
$ws.Cells.Item(6, 12).Value = 6                                         # L6 numero coltelli
$ws.Cells.Item(6, 13).Value = 76                                        # M6 diametro tubo
$ws.Cells.Item(6, 14).Value = 39885                                     # N6 veicolo (now numeric)
$ws.Cells.Item(6, 16).Value = 39885                                     # P6
$ws.Cells.Item(6, 17).Value = "2025-05-09 00:00:00"                     # Q6
$ws.Cells.Item(6, 18).Value = -0.2952660406828704                       # R6
$ws.Cells.Item(6, 19).Value = 1                                         # S6

# --- Row 7 (was commessa 251752, now commessa 251895) ---
$ws.Cells.Item(7, 1).Value  = 251895                                    # A7 commessa
$ws.Cells.Item(7, 3).Value  = 38                                        # C7 minuti setup
$ws.Cells.Item(7, 4).Value  = 249.2112676056338                         # D7 minuti processamento
$ws.Cells.Item(7, 5).Value  = "2025-06-05 07:05:10"                     # E7 inizio setup
$ws.Cells.Item(7, 6).Value  = "2025-06-05 07:43:10"                     # F7 fine setup
$ws.Cells.Item(7, 7).Value  = "2025-06-05 07:43:10"                     # G7 inizio lavorazione
$ws.Cells.Item(7, 8).Value  = "2025-06-05 11:52:23"                     # H7 fine lavorazione
$ws.Cells.Item(7, 9).Value  = 17694                                     # I7 mt da tagliare
$ws.Cells.Item(7, 11).Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9" # K7 macchine compatibili
$ws.Cells.Item(7, 12).Value = 10                                        # L7 numero coltelli
$ws.Cells.Item(7, 13).Value = 70                                        # M7 diametro tubo
$ws.Cells.Item(7, 14).Value = "39891 (esterno)"                         # N7 veicolo (now text)
$ws.Cells.Item(7, 16).Value = 39891                                     # P7
$ws.Cells.Item(7, 17).Value = "2025-05-26 00:00:00"                     # Q7
$ws.Cells.Item(7, 18).Value = -10.49471830986111                        # R7
$ws.Cells.Item(7, 19).Value = 4                                         # S7
